$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: C1 text change
$ws.Range("C1").Value = "class40"

# New data rows (tickers with new "winner" class value of 1)
$data = @(
    @(986,  "MAA",  1),
    @(1103, "TGTX", 1),
    @(1679, "CIB",  1),
    @(1749, "BOH",  1),
    @(1995, "CLSD", 1),
    @(2050, "HASI", 1),
    @(2164, "AVAL", 1),
    @(2171, "SKY",  1),
    @(2383, "NYMX", 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Apply the same formatting as the already-styled ticker-id cell A2 (bold font, thin
# border, centered horizontal, top vertical) to the rest of column A by copying its
# format onto A3:A10.
$ws.Range("A2").Copy()
$ws.Range("A3:A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
